$d = $word.ActiveDocument

# 1) Semester: "Spring 2025" -> "Fall 2025"
$d.Content.Find.Execute("Spring 2025", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Fall 2025", 2) | Out-Null

# 2) Objectives bullet: expand "Review TM4C123 and Keil, learn how to use "
#    into "Review ECE319K, TM4C123, and Keil (or CCS12.8), learn how to use "
$d.Content.Find.Execute("Review TM4C123 and Keil, learn how to use ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Review ECE319K, TM4C123, and Keil (or CCS12.8), learn how to use ", 2) | Out-Null

# 3) Question 3: "Explain how you removed switch bounce" -> "Explain how you removed the switch bounce"
$d.Content.Find.Execute("Explain how you removed switch bounce", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Explain how you removed the switch bounce", 2) | Out-Null

# 4) Question 4 follow-up: "What would have been other approaches?" ->
#    "What are the pros and cons for implementing such a simple software problem with an FSM abstraction?"
$d.Content.Find.Execute("What would have been other approaches?", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "What are the pros and cons for implementing such a simple software problem with an FSM abstraction?", 2) | Out-Null
